# Delete the "노브랜드" row (row 22) from Sheet1, shifting subsequent rows up.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Rows.Item(22).Delete()
